$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 6
$ws.Range("H6").Value = 165.14285
$ws.Range("I6").Value = 165.14285
$ws.Range("K6").Value = 495.42855
$ws.Range("M6").Value = -383.42855

# ALC row 39
$ws.Range("H39").Value = 1417.7333
$ws.Range("I39").Value = 781.4167
$ws.Range("J39").Value = 3963
$ws.Range("K39").Value = 2344.2501
$ws.Range("L39").Value = 11889
$ws.Range("M39").Value = -2048.2501
$ws.Range("N39").Value = -12481

# ALC row 40
$ws.Range("H40").Value = 3049.8333
$ws.Range("I40").Value = 1300
$ws.Range("J40").Value = 3399.8
$ws.Range("K40").Value = 1300
$ws.Range("L40").Value = 3399.8
$ws.Range("M40").Value = -1125
$ws.Range("N40").Value = -3749.8

# ALC row 64
$ws.Range("H64").Value = 7459.25
$ws.Range("I64").Value = 7135.4
$ws.Range("J64").Value = 7999
$ws.Range("K64").Value = 7135.4
$ws.Range("L64").Value = 7999
$ws.Range("M64").Value = -6887.4
$ws.Range("N64").Value = -8495

# ALC row 67
$ws.Range("H67").Value = 7459.25
$ws.Range("I67").Value = 7135.4
$ws.Range("J67").Value = 7999
$ws.Range("K67").Value = 7135.4
$ws.Range("L67").Value = 7999
$ws.Range("M67").Value = -6277.4
$ws.Range("N67").Value = -9715

# ALC row 100
$ws.Range("H100").Value = 12142.75
$ws.Range("I100").Value = 4229.6665
$ws.Range("J100").Value = 16890.6
$ws.Range("K100").Value = 4229.6665
$ws.Range("L100").Value = 16890.6
$ws.Range("M100").Value = -3688.6665
$ws.Range("N100").Value = -17972.6

# ALC row 118
$ws.Range("H118").Value = 935.6667
$ws.Range("I118").Value = 657.0909
$ws.Range("K118").Value = 1971.2727
$ws.Range("M118").Value = -314.2727

# ALC row 129
$ws.Range("H129").Value = 1199.2858
$ws.Range("J129").Value = 2132.6667
$ws.Range("L129").Value = 6398.000100000001
$ws.Range("N129").Value = -16398.0001

# ALC row 132
$ws.Range("H132").Value = 3003.8594
$ws.Range("I132").Value = 2868.3936
$ws.Range("K132").Value = 8605.1808
$ws.Range("M132").Value = -6075.1808

$ws = $wb.Worksheets.Item("ARM")
# ARM row 29
$ws.Range("H29").Value = 2555
$ws.Range("J29").Value = 2555
$ws.Range("L29").Value = 2555
$ws.Range("N29").Value = -3171

# ARM row 61
$ws.Range("H61").Value = 3584.9375
$ws.Range("I61").Value = 3608.4285
$ws.Range("J61").Value = 3566.6667
$ws.Range("K61").Value = 3608.4285
$ws.Range("L61").Value = 3566.6667
$ws.Range("M61").Value = -3396.4285
$ws.Range("N61").Value = -3990.6667

# ARM row 63
$ws.Range("H63").Value = 2796.0908
$ws.Range("I63").Value = 2575.7
$ws.Range("K63").Value = 2575.7
$ws.Range("M63").Value = -1889.7

# ARM row 66
$ws.Range("H66").Value = 2796.0908
$ws.Range("I66").Value = 2575.7
$ws.Range("K66").Value = 12878.5
$ws.Range("M66").Value = -9446.5

# ARM row 74
$ws.Range("H74").Value = 1824.4286
$ws.Range("I74").Value = 1008.86365
$ws.Range("K74").Value = 1008.86365
$ws.Range("M74").Value = -134.86365

# ARM row 77
$ws.Range("H77").Value = 1824.4286
$ws.Range("I77").Value = 1008.86365
$ws.Range("K77").Value = 5044.31825
$ws.Range("M77").Value = -676.3182500000003

# ARM row 122
$ws.Range("H122").Value = 3076.5715
$ws.Range("I122").Value = 1922.8334
$ws.Range("J122").Value = 9999
$ws.Range("K122").Value = 5768.5002
$ws.Range("L122").Value = 29997
$ws.Range("M122").Value = -3318.5002
$ws.Range("N122").Value = -34897

# ARM row 136
$ws.Range("H136").Value = 3584.9375
$ws.Range("I136").Value = 3608.4285
$ws.Range("J136").Value = 3566.6667
$ws.Range("K136").Value = 10825.2855
$ws.Range("L136").Value = 10700.0001
$ws.Range("M136").Value = -8275.2855
$ws.Range("N136").Value = -15800.0001

$ws = $wb.Worksheets.Item("BSM")
# BSM row 12
$ws.Range("H12").Value = 5999
$ws.Range("I12").Value = 5999
$ws.Range("K12").Value = 5999
$ws.Range("M12").Value = -5831

$ws = $wb.Worksheets.Item("CRP")
# CRP row 58
$ws.Range("H58").Value = 1609.9286
$ws.Range("I58").Value = 1314.8334
$ws.Range("K58").Value = 1314.8334
$ws.Range("M58").Value = -1111.8334

# CRP row 132
$ws.Range("H132").Value = 1848.9131
$ws.Range("I132").Value = 2117.9
$ws.Range("J132").Value = 1642
$ws.Range("K132").Value = 6353.700000000001
$ws.Range("L132").Value = 4926
$ws.Range("M132").Value = -3823.700000000001
$ws.Range("N132").Value = -9986

# CRP row 136
$ws.Range("H136").Value = 1609.9286
$ws.Range("I136").Value = 1314.8334
$ws.Range("K136").Value = 3944.5002
$ws.Range("M136").Value = -1394.5002

$ws = $wb.Worksheets.Item("CUL")
# CUL row 37
$ws.Range("H37").Value = 135488.95
$ws.Range("J37").Value = 135488.95
$ws.Range("L37").Value = 406466.85
$ws.Range("N37").Value = -406690.85

# CUL row 122
$ws.Range("H122").Value = 1750.5454
$ws.Range("I122").Value = 212.6
$ws.Range("J122").Value = 3032.1667
$ws.Range("K122").Value = 1913.4
$ws.Range("L122").Value = 27289.5003
$ws.Range("M122").Value = 536.6000000000001
$ws.Range("N122").Value = -32189.5003

# CUL row 131
$ws.Range("H131").Value = 29631658
$ws.Range("J131").Value = 30305128
$ws.Range("L131").Value = 90915384
$ws.Range("N131").Value = -90925464

$ws = $wb.Worksheets.Item("GSM")
# GSM row 80
$ws.Range("H80").Value = 1374.5
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()

# GSM row 83
$ws.Range("H83").Value = 1374.5
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

# GSM row 122
$ws.Range("H122").Value = 4152.6924
$ws.Range("J122").Value = 4806.5386
$ws.Range("L122").Value = 14419.6158
$ws.Range("N122").Value = -19319.6158

# GSM row 132
$ws.Range("H132").Value = 2513.3333
$ws.Range("I132").Value = 2026.6666
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 6079.9998
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -3549.9998
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("LTW")
# LTW row 16
$ws.Range("H16").Value = 557.6842
$ws.Range("I16").Value = 406.2353
$ws.Range("K16").Value = 406.2353
$ws.Range("M16").Value = -236.2353

# LTW row 40
$ws.Range("H40").Value = 2563.4285
$ws.Range("J40").Value = 3049.3333
$ws.Range("L40").Value = 3049.3333
$ws.Range("N40").Value = -3321.3333

# LTW row 68
$ws.Range("H68").Value = 2458.9443
$ws.Range("I68").Value = 2468.5
$ws.Range("J68").Value = 2425.5
$ws.Range("K68").Value = 2468.5
$ws.Range("L68").Value = 2425.5
$ws.Range("M68").Value = -1719.5
$ws.Range("N68").Value = -3923.5

# LTW row 71
$ws.Range("H71").Value = 2458.9443
$ws.Range("I71").Value = 2468.5
$ws.Range("J71").Value = 2425.5
$ws.Range("K71").Value = 12342.5
$ws.Range("L71").Value = 12127.5
$ws.Range("M71").Value = -8598.5
$ws.Range("N71").Value = -19615.5

# LTW row 82
$ws.Range("H82").Value = 3787.3125
$ws.Range("I82").Value = 1685.5
$ws.Range("J82").Value = 4487.9165
$ws.Range("K82").Value = 1685.5
$ws.Range("L82").Value = 4487.9165
$ws.Range("M82").Value = -1324.5
$ws.Range("N82").Value = -5209.9165

# LTW row 85
$ws.Range("H85").Value = 3787.3125
$ws.Range("I85").Value = 1685.5
$ws.Range("J85").Value = 4487.9165
$ws.Range("K85").Value = 1685.5
$ws.Range("L85").Value = 4487.9165
$ws.Range("M85").Value = -437.5
$ws.Range("N85").Value = -6983.9165

# LTW row 93
$ws.Range("H93").Value = 1529
$ws.Range("I93").Value = 615.3333
$ws.Range("J93").Value = 2899.5
$ws.Range("K93").Value = 615.3333
$ws.Range("L93").Value = 2899.5
$ws.Range("M93").Value = 632.6667
$ws.Range("N93").Value = -5395.5

# LTW row 97
$ws.Range("H97").Value = 24424.25
$ws.Range("J97").Value = 24424.25
$ws.Range("L97").Value = 24424.25
$ws.Range("N97").Value = -26406.25

# LTW row 132
$ws.Range("H132").Value = 2278.5527
$ws.Range("J132").Value = 2641.2222
$ws.Range("L132").Value = 7923.6666
$ws.Range("N132").Value = -12983.6666

$ws = $wb.Worksheets.Item("WVR")
# WVR row 62
$ws.Range("H62").Value = 7965.6665
$ws.Range("I62").Value = 7449
$ws.Range("J62").Value = 8999
$ws.Range("K62").Value = 7449
$ws.Range("L62").Value = 8999
$ws.Range("M62").Value = -6825
$ws.Range("N62").Value = -10247

# WVR row 65
$ws.Range("H65").Value = 7965.6665
$ws.Range("I65").Value = 7449
$ws.Range("J65").Value = 8999
$ws.Range("K65").Value = 37245
$ws.Range("L65").Value = 44995
$ws.Range("M65").Value = -34125
$ws.Range("N65").Value = -51235

# WVR row 112
$ws.Range("H112").Value = 23462.334
$ws.Range("J112").Value = 23462.334
$ws.Range("L112").Value = 23462.334
$ws.Range("N112").Value = -26416.334

# WVR row 113
$ws.Range("H113").Value = 2524.6667
$ws.Range("I113").Value = 1161
$ws.Range("K113").Value = 3483
$ws.Range("M113").Value = -1313

# WVR row 117
$ws.Range("H117").Value = 44994
$ws.Range("J117").Value = 44994
$ws.Range("L117").Value = 44994
$ws.Range("N117").Value = -54172

